$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value = "Women's Haircut and Blowout"
$ws.Range("C6").Value = "'$49.99"
$ws.Range("D6").Value = "'$34.99"

$ws.Range("B7").Value = "Gel Manicure"
$ws.Range("C7").Value = "'$34.99"
$ws.Range("D7").Value = "'$29.99"

$ws.Range("C6:D7").ClearFormats()

$ws.Range("F11").Select()
